$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.111.80'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.92%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.460.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.30%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.56'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.70'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.460.23'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.477'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.70%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.405'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.055.25'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.34'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.07%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.462.65'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.150.75'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.46'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.43%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.49'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.608.84'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.183'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.63'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.17'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.00%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.41'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.14'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.39%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '31.98'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +11.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '168.20'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.497.53'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.791'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.39'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.55%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.588.99'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.33'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +11.76%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.63%  '
